$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'63.644.76"
$ws.Range("E2").Value = "  -4.10%  "
$ws.Range("D3").Value = "'3.340.98"
$ws.Range("E3").Value = "  -4.42%  "
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").Value = "'556.32"
$ws.Range("E5").Value = "  -0.10%  "
$ws.Range("D6").Value = "'171.79"
$ws.Range("E6").Value = "  -7.32%  "
$ws.Range("D7").Value = "'0.612"
$ws.Range("E7").Value = "  -4.19%  "
$ws.Range("D8").Value = "'3.332.87"
$ws.Range("E8").Value = "  -4.38%  "
$ws.Range("D9").Value = "'1.00"
$ws.Range("E9").Value = "  -0.09%  "
$ws.Range("D10").Value = "'0.618"
$ws.Range("E10").Value = "  -2.16%  "
$ws.Range("D11").Value = "'0.150"
$ws.Range("E11").Value = "  -1.82%  "
$ws.Range("D12").Value = "'53.69"
$ws.Range("E12").Value = "  -1.41%  "
$ws.Range("D13").Value = "'0.0000265"
$ws.Range("E13").Value = "  -2.27%  "
$ws.Range("D14").Value = "'8.93"
$ws.Range("E14").Value = "  -3.98%  "
$ws.Range("D15").Value = "'3.840.13"
$ws.Range("E15").Value = "  -5.50%  "
$ws.Range("E16").Value = "  -3.29%  "
$ws.Range("D17").Value = "'3.314.39"
$ws.Range("E17").Value = "  -5.27%  "
$ws.Range("D18").Value = "'17.65"
$ws.Range("E18").Value = "  -5.18%  "
$ws.Range("B19").Value = "Uniswap"
$ws.Range("C19").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D19").Value = "'11.63"
$ws.Range("E19").Value = "  -3.14%  "
$ws.Range("B20").Value = "WrappedBTC"
$ws.Range("C20").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D20").Value = "'63.504.84"
$ws.Range("E20").Value = "  -4.40%  "
$ws.Range("D21").Value = "'0.972"
$ws.Range("E21").Value = "  -2.10%  "
$ws.Range("D22").Value = "'404.50"
$ws.Range("E22").Value = "  -4.26%  "
$ws.Range("D23").Value = "'4.07"
$ws.Range("E23").Value = "  +0.28%  "
$ws.Range("D24").Value = "'4.28"
$ws.Range("E24").Value = "  +3.20%  "
$ws.Range("D25").Value = "'82.52"
$ws.Range("E25").Value = "  -4.39%  "
$ws.Range("D26").Value = "'13.06"
$ws.Range("E26").Value = "  +6.22%  "
$ws.Range("D27").Value = "'10.67"
$ws.Range("E27").Value = "  -2.70%  "
$ws.Range("D28").Value = "'2.74"
$ws.Range("E28").Value = "  -5.76%  "
$ws.Range("D29").Value = "'8.74"
$ws.Range("E29").Value = "  -4.02%  "
$ws.Range("D30").Value = "'29.15"
$ws.Range("E30").Value = "  -3.44%  "
$ws.Range("E31").Value = "  -1.57%  "
$ws.Range("D32").Value = "'587.07"
$ws.Range("E32").Value = "  -6.96%  "
$ws.Range("D33").Value = "'11.32"
$ws.Range("E33").Value = "  -3.47%  "
$ws.Range("E34").Value = "  -4.32%  "
$ws.Range("D35").Value = "'57.90"
$ws.Range("E35").Value = "  -3.49%  "
$ws.Range("E36").Value = "  +0.51%  "
$ws.Range("E37").Value = "  +0.17%  "
$ws.Range("D38").Value = "'35.74"
$ws.Range("E38").Value = "  -5.54%  "
$ws.Range("D39").Value = "'3.44"
$ws.Range("E39").Value = "  +0.55%  "
$ws.Range("D40").Value = "'0.0₃0745"
$ws.Range("E40").Value = "  -8.59%  "
$ws.Range("D41").Value = "'0.368"
$ws.Range("E41").Value = "  -4.65%  "
$ws.Range("D42").Value = "'3.135.49"
$ws.Range("E42").Value = "  +0.73%  "
$ws.Range("D44").Value = "'2.83"
$ws.Range("E44").Value = "  -0.70%  "
$ws.Range("E45").Value = "  -4.74%  "
$ws.Range("D46").Value = "'0.0404"
$ws.Range("E46").Value = "  -2.56%  "
$ws.Range("E47").Value = "  -6.23%  "
$ws.Range("D48").Value = "'2.60"
$ws.Range("E48").Value = "  -4.60%  "
$ws.Range("E49").Value = "  -4.13%  "
$ws.Range("D50").Value = "'132.42"
$ws.Range("E50").Value = "  -4.81%  "
$ws.Range("D51").Value = "'8.04"
$ws.Range("E51").Value = "  -5.23%  "
